# ---------------------------------------------------------------------------
# Applies the "subo la funda del antivirus y el plan de antel" commit:
#   1. " Home." -> " Home" + new sentence about the free "Windows defender"
#      antivirus (several runs + a relocated "_GoBack" bookmark).
#   2. Merge the "Impresora: " / "Impresora Canon " runs into a single run.
#   3. Remove the old "_GoBack" bookmark that used to sit at the very end of
#      the document (after "Presupuesto: ... Premium").
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) " Home." -> " Home, y posee el antivirus gratuito<BOOKMARK> "Windows
#    defender", por lo cual el cliente no necesitará costearse un antivirus."
# ---------------------------------------------------------------------------

$marker = [char]0x0001
$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$findText = " Home."
$replacement = " Home, y posee el antivirus gratuito" + $marker + " " + $openQuote + `
    "Windows defender" + $closeQuote + ", por lo cual el cliente no necesitará costearse un antivirus."

$found = $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, `
    $true, 1, $false, $replacement, 2)

if (-not $found) {
    throw "Could not find ' Home.' to replace"
}

# Locate the marker we just inserted and drop a relocated _GoBack bookmark
# right there, then remove the one-character marker itself.
$searchRange = $d.Content
$markerFound = $searchRange.Find.Execute($marker, $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

if (-not $markerFound) {
    throw "Could not find bookmark marker"
}

$bmRange = $d.Range($searchRange.Start, $searchRange.Start)

# Remove the stale _GoBack bookmark (currently at the end of the document)
# before re-adding it at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($searchRange.Start, $searchRange.Start + 1)
$markerRange.Delete()

# ---------------------------------------------------------------------------
# 2) "Impresora: " + "Impresora Canon " -> "Impresora: Impresora Canon "
#    (merge the two adjacent, identically-formatted runs into one)
# ---------------------------------------------------------------------------

$printerFound = $d.Content.Find.Execute("Impresora: Impresora Canon ", $false, $false, `
    $false, $false, $false, $true, 1, $false, "Impresora: Impresora Canon ", 2)

if (-not $printerFound) {
    throw "Could not find 'Impresora: Impresora Canon ' text"
}
